$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily log. It belongs right
# after the current row 17, so insert a fresh row at 18 (pushing the old
# rows 18-57 down to 19-58) and populate it with the new observation.
$ws.Rows("18").Insert()

$ws.Range("A18").Value = 10
$ws.Range("B18").Value = 'Vega Modelo de Temuco'
$ws.Range("C18").Value = 'La Araucanía'
$ws.Range("D18").Value = 44930
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112042
$ws.Range("G18").Value = 'Locoto'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2500
$ws.Range("N18").Value = '$/kilo'
$ws.Range("O18").Value = 'Región de Arica y Parinacota'
$ws.Range("P18").Value = 2500
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 'Hortaliza'
